$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.81824933333333
$ws.Range("H2").Value = 53.454748
$ws.Range("I2").Value = 0.05180179233147032
$ws.Range("J2").Value = 0.05180179233147034
$ws.Range("M2").Value = 1.471908333333333
$ws.Range("N2").Value = 4.415725
$ws.Range("O2").Value = 0.2507683239665115
$ws.Range("P2").Value = 0.2507683239665115
$ws.Range("Q2").Value = 26.22682967914445
$ws.Range("R2").Value = 236.0414671123
$ws.Range("S2").Value = 0.0129902486414241
$ws.Range("T2").Value = 0.01299024864142411

$ws.Range("G3").Value = 17.81824933333333
$ws.Range("H3").Value = 53.454748
$ws.Range("I3").Value = 0.05180179233147032
$ws.Range("J3").Value = 0.05180179233147034
$ws.Range("N3").Value = 7.040756999999999
$ws.Range("O3").Value = 0.3998434758381655
$ws.Range("P3").Value = 0.3998434758381655
$ws.Range("Q3").Value = 41.81798790713733
$ws.Range("R3").Value = 376.361891164236
$ws.Range("S3").Value = 0.02071260870046192
$ws.Range("T3").Value = 0.02071260870046193

$ws.Range("G4").Value = 17.81824933333333
$ws.Range("H4").Value = 53.454748
$ws.Range("I4").Value = 0.05180179233147032
$ws.Range("J4").Value = 0.05180179233147034
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.8927443333333333
$ws.Range("N4").Value = 2.678233
$ws.Range("O4").Value = 0.152096428242656
$ws.Range("P4").Value = 0.1520964282426559
$ws.Range("Q4").Value = 15.90714112225378
$ws.Range("R4").Value = 143.164270100284
$ws.Range("S4").Value = 0.007878867590184441
$ws.Range("T4").Value = 0.007878867590184442

$ws.Range("G5").Value = 17.81824933333333
$ws.Range("H5").Value = 53.454748
$ws.Range("I5").Value = 0.05180179233147032
$ws.Range("J5").Value = 0.05180179233147034
$ws.Range("M5").Value = 1.158022666666667
$ws.Range("N5").Value = 3.474068
$ws.Range("O5").Value = 0.1972917719526671
$ws.Range("P5").Value = 0.197291771952667
$ws.Range("Q5").Value = 20.63393660831822
$ws.Range("R5").Value = 185.705429474864
$ws.Range("S5").Value = 0.01022006739939986
$ws.Range("T5").Value = 0.01022006739939986

$ws.Range("I6").Value = 0.4402211587141748
$ws.Range("J6").Value = 0.4402211587141748
$ws.Range("M6").Value = 1.471908333333333
$ws.Range("N6").Value = 4.415725
$ws.Range("O6").Value = 0.2507683239665115
$ws.Range("P6").Value = 0.2507683239665115
$ws.Range("Q6").Value = 222.8804223003334
$ws.Range("R6").Value = 2005.923800703
$ws.Range("S6").Value = 0.1103935221453493
$ws.Range("T6").Value = 0.1103935221453493

$ws.Range("I7").Value = 0.4402211587141748
$ws.Range("J7").Value = 0.4402211587141748
$ws.Range("N7").Value = 7.040756999999999
$ws.Range("O7").Value = 0.3998434758381655
$ws.Range("P7").Value = 0.3998434758381655
$ws.Range("Q7").Value = 355.37695247644
$ws.Range("S7").Value = 0.1760195582377804
$ws.Range("T7").Value = 0.1760195582377804

$ws.Range("I8").Value = 0.4402211587141748
$ws.Range("J8").Value = 0.4402211587141748
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.8927443333333333
$ws.Range("N8").Value = 2.678233
$ws.Range("O8").Value = 0.152096428242656
$ws.Range("P8").Value = 0.1520964282426559
$ws.Range("Q8").Value = 135.1818109276933
$ws.Range("R8").Value = 1216.63629834924
$ws.Range("S8").Value = 0.06695606587726935
$ws.Range("T8").Value = 0.06695606587726934

$ws.Range("I9").Value = 0.4402211587141748
$ws.Range("J9").Value = 0.4402211587141748
$ws.Range("M9").Value = 1.158022666666667
$ws.Range("N9").Value = 3.474068
$ws.Range("O9").Value = 0.1972917719526671
$ws.Range("P9").Value = 0.197291771952667
$ws.Range("Q9").Value = 175.3509883292267
$ws.Range("R9").Value = 1578.15889496304
$ws.Range("S9").Value = 0.08685201245377583
$ws.Range("T9").Value = 0.08685201245377581

$ws.Range("G10").Value = 100.6958183333333
$ws.Range("H10").Value = 302.087455
$ws.Range("I10").Value = 0.29274614875843
$ws.Range("J10").Value = 0.2927461487584301
$ws.Range("M10").Value = 1.471908333333333
$ws.Range("N10").Value = 4.415725
$ws.Range("O10").Value = 0.2507683239665115
$ws.Range("P10").Value = 0.2507683239665115
$ws.Range("Q10").Value = 148.2150141366528
$ws.Range("R10").Value = 1333.935127229875
$ws.Range("S10").Value = 0.07341146107180255
$ws.Range("T10").Value = 0.07341146107180256

$ws.Range("G11").Value = 100.6958183333333
$ws.Range("H11").Value = 302.087455
$ws.Range("I11").Value = 0.29274614875843
$ws.Range("J11").Value = 0.2927461487584301
$ws.Range("N11").Value = 7.040756999999999
$ws.Range("O11").Value = 0.3998434758381655
$ws.Range("P11").Value = 0.3998434758381655
$ws.Range("Q11").Value = 236.3249292670483
$ws.Range("R11").Value = 2126.924363403435
$ws.Range("S11").Value = 0.1170526376578073
$ws.Range("T11").Value = 0.1170526376578074

$ws.Range("G12").Value = 100.6958183333333
$ws.Range("H12").Value = 302.087455
$ws.Range("I12").Value = 0.29274614875843
$ws.Range("J12").Value = 0.2927461487584301
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.8927443333333333
$ws.Range("N12").Value = 2.678233
$ws.Range("O12").Value = 0.152096428242656
$ws.Range("P12").Value = 0.1520964282426559
$ws.Range("Q12").Value = 89.89562120744611
$ws.Range("R12").Value = 809.060590867015
$ws.Range("S12").Value = 0.04452564360795044
$ws.Range("T12").Value = 0.04452564360795044

$ws.Range("G13").Value = 100.6958183333333
$ws.Range("H13").Value = 302.087455
$ws.Range("I13").Value = 0.29274614875843
$ws.Range("J13").Value = 0.2927461487584301
$ws.Range("M13").Value = 1.158022666666667
$ws.Range("N13").Value = 3.474068
$ws.Range("O13").Value = 0.1972917719526671
$ws.Range("P13").Value = 0.197291771952667
$ws.Range("Q13").Value = 116.6080400685489
$ws.Range("R13").Value = 1049.47236061694
$ws.Range("S13").Value = 0.05775640642086972
$ws.Range("T13").Value = 0.05775640642086972

$ws.Range("G14").Value = 74.032918
$ws.Range("H14").Value = 222.098754
$ws.Range("I14").Value = 0.2152309001959248
$ws.Range("J14").Value = 0.2152309001959249
$ws.Range("M14").Value = 1.471908333333333
$ws.Range("N14").Value = 4.415725
$ws.Range("O14").Value = 0.2507683239665115
$ws.Range("P14").Value = 0.2507683239665115
$ws.Range("Q14").Value = 108.9696689451833
$ws.Range("R14").Value = 980.72702050665
$ws.Range("S14").Value = 0.05397309210793558
$ws.Range("T14").Value = 0.0539730921079356

$ws.Range("G15").Value = 74.032918
$ws.Range("H15").Value = 222.098754
$ws.Range("I15").Value = 0.2152309001959248
$ws.Range("J15").Value = 0.2152309001959249
$ws.Range("N15").Value = 7.040756999999999
$ws.Range("O15").Value = 0.3998434758381655
$ws.Range("P15").Value = 0.3998434758381655
$ws.Range("Q15").Value = 173.749261879642
$ws.Range("R15").Value = 1563.743356916778
$ws.Range("S15").Value = 0.08605867124211589
$ws.Range("T15").Value = 0.08605867124211591

$ws.Range("G16").Value = 74.032918
$ws.Range("H16").Value = 222.098754
$ws.Range("I16").Value = 0.2152309001959248
$ws.Range("J16").Value = 0.2152309001959249
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.8927443333333333
$ws.Range("N16").Value = 2.678233
$ws.Range("O16").Value = 0.152096428242656
$ws.Range("P16").Value = 0.1520964282426559
$ws.Range("Q16").Value = 66.09246802463132
$ws.Range("R16").Value = 594.8322122216819
$ws.Range("S16").Value = 0.03273585116725173
$ws.Range("T16").Value = 0.03273585116725173

$ws.Range("G17").Value = 74.032918
$ws.Range("H17").Value = 222.098754
$ws.Range("I17").Value = 0.2152309001959248
$ws.Range("J17").Value = 0.2152309001959249
$ws.Range("M17").Value = 1.158022666666667
$ws.Range("N17").Value = 3.474068
$ws.Range("O17").Value = 0.1972917719526671
$ws.Range("P17").Value = 0.197291771952667
$ws.Range("Q17").Value = 85.73179712347466
$ws.Range("R17").Value = 771.5861741112719
$ws.Range("S17").Value = 0.04246328567862165
$ws.Range("T17").Value = 0.04246328567862165
